$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 'This graph depicts the TKN_About in TKN_Geo in TKN_Year TKN_UOM meters.'
$ws.Range("X5").Value = 'There is another high of about TKN_UOM meters in November.'
$ws.Range("X6").Value = 'This graph depicts the TKN_About in TKN_UOM of TKN_Geo in TKN_Year.'
$ws.Range("X7").Value = 'Production starts at approximately 64 TKN_UOM in January and increases to one of its maximum values of 102 TKN_UOM in March.'
$ws.Range("X9").Value = 'From July until TKN_About gradually increases again until it reaches its second maximum of 102 TKN_UOM in November.'
$ws.Range("X10").Value = 'There is a sharp drop from November onwards, reaching its minimum value of -1 TKN_UOM in December.'
$ws.Range("X11").Value = 'This line graph displays TKN_About in TKN_UOM in TKN_Geo for TKN_Year.'
$ws.Range("X14").Value = 'A steady increase from August to November is followed by a dramatic reduction in December to a yearly low of TKN_UOM metres.'
$ws.Range("X15").Value = 'The following graph shows information about TKN_About TKN_UOM in TKN_Geo during TKN_Year.'
$ws.Range("X16").Value = 'As can be seen from the graph, TKN_UOM have been steadily increasing its value over the year.'
$ws.Range("X18").Value = 'This graph represents TKN_About TKN_UOM of TKN_Geo in TKN_Year.'
$ws.Range("X20").Value = 'Continuous increase of TKN_About TKN_UOM in TKN_Geo during the year TKN_Year.'
$ws.Range("X21").Value = 'The graph shows TKN_UOM in TKN_Geo for the year TKN_Year.'
$ws.Range("X24").Value = 'The TKN_About from the year TKN_Year is showen by the graph.'
$ws.Range("X28").Value = 'The shown data are about TKN_UOM in TKN_Geo during TKN_Year.'
$ws.Range("X32").Value = 'Line chart showing TKN_About TKN_UOM in TKN_Geo in TKN_Year.'
$ws.Range("X38").Value = 'The graph shows the TKN_About fruits in TKN_Geo over the year of TKN_Year.'
$ws.Range("X39").Value = 'In the first half TKN_About was nearly constant, until it fall dramatically in the mounth of August.'
$ws.Range("X40").Value = 'Afterwards the TKN_About rises again over the index of 75.'
$ws.Range("X41").Value = 'The graph is a line chart, illustrating the TKN_About in TKN_Geo.'
$ws.Range("X48").Value = 'The TKN_About in TKN_Geo considerably increased over TKN_Year.'
$ws.Range("X52").Value = 'It shows the TKN_UOM of TKN_About TKN_Geo during TKN_Year.'
$ws.Range("X54").Value = 'After that, the TKN_UOM TKN_About TKN_Geo decrease during the last few months of the year.'
$ws.Range("X56").Value = 'The graph showes the TKN_About in TKN_Year.'
$ws.Range("X59").Value = 'From August till December the graph decrease continuslie TKN_About TKN_UOM of 28 in December.'
$ws.Range("X63").Value = 'The line chart describes the TKN_UOM of TKN_About TKN_Geo during TKN_Year.'
$ws.Range("X65").Value = 'During the last few months of the year the TKN_UOM of TKN_About TKN_Geo rapidly decreased.'
$ws.Range("X66").Value = 'The TKN_UOM TKN_About TKN_Geo during TKN_Year substantially increased over the summer months.'
$ws.Range("X69").Value = 'The graph illustrates the TKN_About in TKN_Geo for the year TKN_Year.'
$ws.Range("X70").Value = 'There are sharp decreases in TKN_About February, July and November.'
$ws.Range("X71").Value = 'This graph shows the monthly amount TKN_About produced in TKN_Geo in TKN_Year in tonnes.'
$ws.Range("X72").Value = 'Oat TKN_About 0 and 97 TKN_UOM each month.'
$ws.Range("X73").Value = 'While most months TKN_About at around 45 TKN_UOM, February, June and July saw dips in production to around 7 TKN_UOM and the fall months saw a higher harvest, -48 and above 90 TKN_UOM consecutively.'
$ws.Range("X74").Value = 'The following graph depicts the Canadian TKN_About in TKN_UOM during TKN_Year.'
$ws.Range("X75").Value = 'In January approximately 60 TKN_UOM were produced, after which a sharp drop can be observed for February, reaching slightly below 15.'
$ws.Range("X78").Value = 'Following this there is a drop, TKN_About remains almost unchanged during November and December at 45.'
$ws.Range("X79").Value = 'This graph depicts the TKN_About in TKN_UOM in TKN_Geo during TKN_Year.'
$ws.Range("X81").Value = 'For the following TKN_About slightly increases until August.'
$ws.Range("X82").Value = 'Between August and September a rapid increase TKN_About can be observed, indicating the maximum of over 96 in September.'
$ws.Range("X84").Value = 'The TKN_About strongly fluctuated over the TKN_Year in TKN_Geo.'
$ws.Range("X88").Value = 'The figure shows a line chart about the TKN_About (in tonnes) in TKN_Geo during TKN_Year.'
$ws.Range("X90").Value = 'Even that, if it has been recorded a significant growth from November to December, where the TKN_About during TKN_Year reached its peak.'
$ws.Range("X91").Value = 'The graph shows the TKN_About in TKN_UOM in TKN_Geo in TKN_Year.'
$ws.Range("X97").Value = 'The line chart is displaying the recorded data about TKN_About in TKN_Geo during TKN_Year.'
$ws.Range("X98").Value = 'The minimum TKN_About have been recorded during the months of January, August and November.'
$ws.Range("X100").Value = 'The maximum TKN_About over the year was recorded during December.'
$ws.Range("X101").Value = 'This graph represents the TKN_About in TKN_Geo in TKN_Year (in tonnes).'
$ws.Range("X103").Value = 'There is a sharp increase during February and during the TKN_About levels off between 42 and 61.'
$ws.Range("X105").Value = 'For the following -116 months TKN_About remains steady at approximately 47.'
$ws.Range("X106").Value = 'After a drop TKN_About peaks in December with up to 101 tonnes.'
